# Auto-generated Excel COM-interop script to apply scheduled market-price updates
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 48362.09  # H28: 48363.227 -> 48362.09
$ws.Cells.Item(28, 9).Value = 63658.812  # I28: 67889.07000000001 -> 63658.812
$ws.Cells.Item(28, 10).Value = 7570.8335  # J28: 6522.143 -> 7570.8335
$ws.Cells.Item(28, 11).Value = 63658.812  # K28: 67889.07000000001 -> 63658.812
$ws.Cells.Item(28, 12).Value = 7570.8335  # L28: 6522.143 -> 7570.8335
$ws.Cells.Item(28, 13).Value = -63173.812  # M28: -67404.07000000001 -> -63173.812
$ws.Cells.Item(28, 14).Value = -8540.833500000001  # N28: -7492.143 -> -8540.833500000001
$ws.Cells.Item(106, 8).Value = 4755  # H106: 4913.85 -> 4755
$ws.Cells.Item(106, 9).Value = 2367.2942  # I106: 2416.625 -> 2367.2942
$ws.Cells.Item(106, 11).Value = 2367.2942  # K106: 2416.625 -> 2367.2942
$ws.Cells.Item(106, 13).Value = -1736.2942  # M106: -1785.625 -> -1736.2942
$ws.Cells.Item(129, 8).Value = 937.1667  # H129: 1581.4 -> 937.1667
$ws.Cells.Item(129, 9).Value = 937.1667  # I129: 1402 -> 937.1667
$ws.Cells.Item(129, 10).Value = 0  # J129: 2000 -> 0
$ws.Cells.Item(129, 11).Value = 2811.5001  # K129: 4206 -> 2811.5001
$ws.Cells.Item(129, 12).Value = 0  # L129: 6000 -> 0
$ws.Cells.Item(129, 13).Value = 2188.4999  # M129: 794 -> 2188.4999
$ws.Cells.Item(129, 14).ClearContents()  # N129: -16000 -> (removed)
$ws.Cells.Item(137, 8).Value = 2626.617  # H137: 2424.8064 -> 2626.617
$ws.Cells.Item(137, 9).Value = 2920.2  # I137: 1966.0667 -> 2920.2
$ws.Cells.Item(137, 10).Value = 2591.6667  # J137: 2571.2126 -> 2591.6667
$ws.Cells.Item(137, 11).Value = 8760.599999999999  # K137: 5898.2001 -> 8760.599999999999
$ws.Cells.Item(137, 12).Value = 7775.000100000001  # L137: 7713.6378 -> 7775.000100000001
$ws.Cells.Item(137, 13).Value = -6210.599999999999  # M137: -3348.2001 -> -6210.599999999999
$ws.Cells.Item(137, 14).Value = -12875.0001  # N137: -12813.6378 -> -12875.0001
$ws.Cells.Item(138, 8).Value = 3192.9395  # H138: 2453.255 -> 3192.9395
$ws.Cells.Item(138, 9).Value = 3869.25  # I138: 1675.24 -> 3869.25
$ws.Cells.Item(138, 10).Value = 3099.6553  # J138: 3201.3462 -> 3099.6553
$ws.Cells.Item(138, 11).Value = 11607.75  # K138: 5025.72 -> 11607.75
$ws.Cells.Item(138, 12).Value = 9298.965899999999  # L138: 9604.0386 -> 9298.965899999999
$ws.Cells.Item(138, 13).Value = -6467.75  # M138: 114.2799999999997 -> -6467.75
$ws.Cells.Item(138, 14).Value = -19578.9659  # N138: -19884.0386 -> -19578.9659
$ws.Cells.Item(141, 8).Value = 2500.25  # H141: 2180.1428 -> 2500.25
$ws.Cells.Item(141, 9).Value = 2444.7778  # I141: 2117.1155 -> 2444.7778
$ws.Cells.Item(141, 11).Value = 7334.3334  # K141: 6351.3465 -> 7334.3334
$ws.Cells.Item(141, 13).Value = -2154.3334  # M141: -1171.3465 -> -2154.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1818.8235  # H2: 1952.0714 -> 1818.8235
$ws.Cells.Item(2, 9).Value = 1718.3636  # I2: 1913.875 -> 1718.3636
$ws.Cells.Item(2, 11).Value = 1718.3636  # K2: 1913.875 -> 1718.3636
$ws.Cells.Item(2, 13).Value = -1605.3636  # M2: -1800.875 -> -1605.3636
$ws.Cells.Item(32, 8).Value = 109128.83  # H32: 138673.94 -> 109128.83
$ws.Cells.Item(32, 9).Value = 119124.21  # I32: 155319.81 -> 119124.21
$ws.Cells.Item(32, 10).Value = 2927.875  # J32: 3426.125 -> 2927.875
$ws.Cells.Item(32, 11).Value = 119124.21  # K32: 155319.81 -> 119124.21
$ws.Cells.Item(32, 12).Value = 2927.875  # L32: 3426.125 -> 2927.875
$ws.Cells.Item(32, 13).Value = -118837.21  # M32: -155032.81 -> -118837.21
$ws.Cells.Item(32, 14).Value = -3501.875  # N32: -4000.125 -> -3501.875
$ws.Cells.Item(43, 8).Value = 25000  # H43: 0 -> 25000
$ws.Cells.Item(43, 10).Value = 25000  # J43: 0 -> 25000
$ws.Cells.Item(43, 12).Value = 25000  # L43: 0 -> 25000
$ws.Cells.Item(43, 14).Value = -25626  # N43: None -> -25626
$ws.Cells.Item(57, 8).Value = 21739130  # H57: 10874565 -> 21739130
$ws.Cells.Item(57, 9).Value = 21739130  # I57: 10874565 -> 21739130
$ws.Cells.Item(57, 11).Value = 21739130  # K57: 10874565 -> 21739130
$ws.Cells.Item(57, 13).Value = -21738646  # M57: -10874081 -> -21738646
$ws.Cells.Item(61, 8).Value = 2402.5334  # H61: 1002717.4 -> 2402.5334
$ws.Cells.Item(61, 9).Value = 1594.409  # I61: 834945.7 -> 1594.409
$ws.Cells.Item(61, 10).Value = 4624.875  # J61: 1254374.9 -> 4624.875
$ws.Cells.Item(61, 11).Value = 1594.409  # K61: 834945.7 -> 1594.409
$ws.Cells.Item(61, 12).Value = 4624.875  # L61: 1254374.9 -> 4624.875
$ws.Cells.Item(61, 13).Value = -1382.409  # M61: -834733.7 -> -1382.409
$ws.Cells.Item(61, 14).Value = -5048.875  # N61: -1254798.9 -> -5048.875
$ws.Cells.Item(110, 8).Value = 35715440  # H110: 50001070 -> 35715440
$ws.Cells.Item(110, 9).Value = 41667776  # I110: 55556580 -> 41667776
$ws.Cells.Item(110, 10).Value = 1427  # J110: 1506 -> 1427
$ws.Cells.Item(110, 11).Value = 41667776  # K110: 55556580 -> 41667776
$ws.Cells.Item(110, 12).Value = 1427  # L110: 1506 -> 1427
$ws.Cells.Item(110, 13).Value = -41665731  # M110: -55554535 -> -41665731
$ws.Cells.Item(110, 14).Value = -5517  # N110: -5596 -> -5517
$ws.Cells.Item(116, 8).Value = 1818.8235  # H116: 1952.0714 -> 1818.8235
$ws.Cells.Item(116, 9).Value = 1718.3636  # I116: 1913.875 -> 1718.3636
$ws.Cells.Item(116, 11).Value = 1718.3636  # K116: 1913.875 -> 1718.3636
$ws.Cells.Item(116, 13).Value = 575.6364000000001  # M116: 380.125 -> 575.6364000000001
$ws.Cells.Item(122, 8).Value = 2569.6086  # H122: 2193.7 -> 2569.6086
$ws.Cells.Item(122, 9).Value = 1711.8235  # I122: 2193.7 -> 1711.8235
$ws.Cells.Item(122, 10).Value = 5000  # J122: 0 -> 5000
$ws.Cells.Item(122, 11).Value = 5135.470499999999  # K122: 6581.099999999999 -> 5135.470499999999
$ws.Cells.Item(122, 12).Value = 15000  # L122: 0 -> 15000
$ws.Cells.Item(122, 13).Value = -2685.470499999999  # M122: -4131.099999999999 -> -2685.470499999999
$ws.Cells.Item(122, 14).Value = -19900  # N122: None -> -19900
$ws.Cells.Item(132, 8).Value = 489150.8  # H132: 466441.12 -> 489150.8
$ws.Cells.Item(132, 9).Value = 286919.53  # I132: 278957.12 -> 286919.53
$ws.Cells.Item(132, 10).Value = 1668833.1  # J132: 1430644.8 -> 1668833.1
$ws.Cells.Item(132, 11).Value = 860758.5900000001  # K132: 836871.36 -> 860758.5900000001
$ws.Cells.Item(132, 12).Value = 5006499.300000001  # L132: 4291934.4 -> 5006499.300000001
$ws.Cells.Item(132, 13).Value = -858228.5900000001  # M132: -834341.36 -> -858228.5900000001
$ws.Cells.Item(132, 14).Value = -5011559.300000001  # N132: -4296994.4 -> -5011559.300000001
$ws.Cells.Item(136, 8).Value = 2402.5334  # H136: 1002717.4 -> 2402.5334
$ws.Cells.Item(136, 9).Value = 1594.409  # I136: 834945.7 -> 1594.409
$ws.Cells.Item(136, 10).Value = 4624.875  # J136: 1254374.9 -> 4624.875
$ws.Cells.Item(136, 11).Value = 4783.227000000001  # K136: 2504837.1 -> 4783.227000000001
$ws.Cells.Item(136, 12).Value = 13874.625  # L136: 3763124.7 -> 13874.625
$ws.Cells.Item(136, 13).Value = -2233.227000000001  # M136: -2502287.1 -> -2233.227000000001
$ws.Cells.Item(136, 14).Value = -18974.625  # N136: -3768224.7 -> -18974.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1818.8235  # H3: 1952.0714 -> 1818.8235
$ws.Cells.Item(3, 9).Value = 1718.3636  # I3: 1913.875 -> 1718.3636
$ws.Cells.Item(3, 11).Value = 1718.3636  # K3: 1913.875 -> 1718.3636
$ws.Cells.Item(3, 13).Value = -1604.3636  # M3: -1799.875 -> -1604.3636
$ws.Cells.Item(20, 8).Value = 2885.4285  # H20: 2875.8096 -> 2885.4285
$ws.Cells.Item(20, 9).Value = 2817.5  # I20: 2702.625 -> 2817.5
$ws.Cells.Item(20, 10).Value = 3021.2856  # J20: 3430 -> 3021.2856
$ws.Cells.Item(20, 11).Value = 2817.5  # K20: 2702.625 -> 2817.5
$ws.Cells.Item(20, 12).Value = 3021.2856  # L20: 3430 -> 3021.2856
$ws.Cells.Item(20, 13).Value = -2570.5  # M20: -2455.625 -> -2570.5
$ws.Cells.Item(20, 14).Value = -3515.2856  # N20: -3924 -> -3515.2856
$ws.Cells.Item(105, 8).Value = 2262.8235  # H105: 1920.7307 -> 2262.8235
$ws.Cells.Item(105, 9).Value = 2155.5557  # I105: 1664 -> 2155.5557
$ws.Cells.Item(105, 10).Value = 2383.5  # J105: 2405.6667 -> 2383.5
$ws.Cells.Item(105, 11).Value = 2155.5557  # K105: 1664 -> 2155.5557
$ws.Cells.Item(105, 12).Value = 2383.5  # L105: 2405.6667 -> 2383.5
$ws.Cells.Item(105, 13).Value = -408.5556999999999  # M105: 83 -> -408.5556999999999
$ws.Cells.Item(105, 14).Value = -5877.5  # N105: -5899.6667 -> -5877.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 603.5  # H7: 398.4375 -> 603.5
$ws.Cells.Item(7, 9).Value = 524.2  # I7: 407 -> 524.2
$ws.Cells.Item(7, 10).Value = 1000  # J7: 361.33334 -> 1000
$ws.Cells.Item(7, 11).Value = 524.2  # K7: 407 -> 524.2
$ws.Cells.Item(7, 12).Value = 1000  # L7: 361.33334 -> 1000
$ws.Cells.Item(7, 13).Value = -411.2  # M7: -294 -> -411.2
$ws.Cells.Item(7, 14).Value = -1226  # N7: -587.33334 -> -1226
$ws.Cells.Item(31, 8).Value = 2651.6555  # H31: 2561.8765 -> 2651.6555
$ws.Cells.Item(31, 9).Value = 1573.8235  # I31: 1734.9286 -> 1573.8235
$ws.Cells.Item(31, 10).Value = 2902.6575  # J31: 2734.6716 -> 2902.6575
$ws.Cells.Item(31, 11).Value = 1573.8235  # K31: 1734.9286 -> 1573.8235
$ws.Cells.Item(31, 12).Value = 2902.6575  # L31: 2734.6716 -> 2902.6575
$ws.Cells.Item(31, 13).Value = -1278.8235  # M31: -1439.9286 -> -1278.8235
$ws.Cells.Item(31, 14).Value = -3492.6575  # N31: -3324.6716 -> -3492.6575
$ws.Cells.Item(34, 8).Value = 2651.6555  # H34: 2561.8765 -> 2651.6555
$ws.Cells.Item(34, 9).Value = 1573.8235  # I34: 1734.9286 -> 1573.8235
$ws.Cells.Item(34, 10).Value = 2902.6575  # J34: 2734.6716 -> 2902.6575
$ws.Cells.Item(34, 11).Value = 1573.8235  # K34: 1734.9286 -> 1573.8235
$ws.Cells.Item(34, 12).Value = 2902.6575  # L34: 2734.6716 -> 2902.6575
$ws.Cells.Item(34, 13).Value = -1371.8235  # M34: -1532.9286 -> -1371.8235
$ws.Cells.Item(34, 14).Value = -3306.6575  # N34: -3138.6716 -> -3306.6575
$ws.Cells.Item(93, 8).Value = 0  # H93: 4000 -> 0
$ws.Cells.Item(93, 9).Value = 0  # I93: 4000 -> 0
$ws.Cells.Item(93, 11).Value = 0  # K93: 4000 -> 0
$ws.Cells.Item(93, 13).ClearContents()  # M93: -2128 -> (removed)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 44992944  # H4: 56324500 -> 44992944
$ws.Cells.Item(4, 9).Value = 44992944  # I4: 56324500 -> 44992944
$ws.Cells.Item(4, 11).Value = 134978832  # K4: 168973500 -> 134978832
$ws.Cells.Item(4, 13).Value = -134978720  # M4: -168973388 -> -134978720
$ws.Cells.Item(8, 8).Value = 4316.2144  # H8: 4189.3125 -> 4316.2144
$ws.Cells.Item(8, 9).Value = 4316.2144  # I8: 4189.3125 -> 4316.2144
$ws.Cells.Item(8, 11).Value = 12948.6432  # K8: 12567.9375 -> 12948.6432
$ws.Cells.Item(8, 13).Value = -12809.6432  # M8: -12428.9375 -> -12809.6432
$ws.Cells.Item(11, 8).Value = 38375  # H11: 79817.66 -> 38375
$ws.Cells.Item(11, 9).Value = 54921.58  # I11: 69557.87 -> 54921.58
$ws.Cells.Item(11, 10).Value = 3443.3333  # J11: 93808.27 -> 3443.3333
$ws.Cells.Item(11, 11).Value = 164764.74  # K11: 208673.61 -> 164764.74
$ws.Cells.Item(11, 12).Value = 10329.9999  # L11: 281424.81 -> 10329.9999
$ws.Cells.Item(11, 13).Value = -164624.74  # M11: -208533.61 -> -164624.74
$ws.Cells.Item(11, 14).Value = -10609.9999  # N11: -281704.81 -> -10609.9999
$ws.Cells.Item(14, 8).Value = 681.0769  # H14: 734 -> 681.0769
$ws.Cells.Item(14, 9).Value = 681.0769  # I14: 734 -> 681.0769
$ws.Cells.Item(14, 11).Value = 2043.2307  # K14: 2202 -> 2043.2307
$ws.Cells.Item(14, 13).Value = -1870.2307  # M14: -2029 -> -1870.2307
$ws.Cells.Item(98, 8).Value = 5561  # H98: 5617.4375 -> 5561
$ws.Cells.Item(98, 9).Value = 9346.429  # I98: 9698.429 -> 9346.429
$ws.Cells.Item(98, 10).Value = 2248.75  # J98: 2443.3333 -> 2248.75
$ws.Cells.Item(98, 11).Value = 28039.287  # K98: 29095.287 -> 28039.287
$ws.Cells.Item(98, 12).Value = 6746.25  # L98: 7329.999899999999 -> 6746.25
$ws.Cells.Item(98, 13).Value = -26541.287  # M98: -27597.287 -> -26541.287
$ws.Cells.Item(98, 14).Value = -9742.25  # N98: -10325.9999 -> -9742.25
$ws.Cells.Item(132, 8).Value = 999.5  # H132: 0 -> 999.5
$ws.Cells.Item(132, 10).Value = 999.5  # J132: 0 -> 999.5
$ws.Cells.Item(132, 12).Value = 8995.5  # L132: 0 -> 8995.5
$ws.Cells.Item(132, 14).Value = -14055.5  # N132: None -> -14055.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2693.625  # H102: 2396.0667 -> 2693.625
$ws.Cells.Item(102, 9).Value = 1940.4375  # I102: 2278.4167 -> 1940.4375
$ws.Cells.Item(102, 10).Value = 4200  # J102: 2866.6667 -> 4200
$ws.Cells.Item(102, 11).Value = 1940.4375  # K102: 2278.4167 -> 1940.4375
$ws.Cells.Item(102, 12).Value = 4200  # L102: 2866.6667 -> 4200
$ws.Cells.Item(102, 13).Value = -318.4375  # M102: -656.4167000000002 -> -318.4375
$ws.Cells.Item(102, 14).Value = -7444  # N102: -6110.6667 -> -7444
$ws.Cells.Item(122, 8).Value = 2671.2415  # H122: 3224.4211 -> 2671.2415
$ws.Cells.Item(122, 9).Value = 2515.4546  # I122: 2947.4285 -> 2515.4546
$ws.Cells.Item(122, 10).Value = 3160.8572  # J122: 4000 -> 3160.8572
$ws.Cells.Item(122, 11).Value = 7546.3638  # K122: 8842.2855 -> 7546.3638
$ws.Cells.Item(122, 12).Value = 9482.571599999999  # L122: 12000 -> 9482.571599999999
$ws.Cells.Item(122, 13).Value = -5096.3638  # M122: -6392.2855 -> -5096.3638
$ws.Cells.Item(122, 14).Value = -14382.5716  # N122: -16900 -> -14382.5716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 8763.409  # H61: 425026.75 -> 8763.409
$ws.Cells.Item(61, 9).Value = 9713.333000000001  # I61: 509134.34 -> 9713.333000000001
$ws.Cells.Item(61, 11).Value = 9713.333000000001  # K61: 509134.34 -> 9713.333000000001
$ws.Cells.Item(61, 13).Value = -9511.333000000001  # M61: -508932.34 -> -9511.333000000001
$ws.Cells.Item(113, 8).Value = 8763.409  # H113: 425026.75 -> 8763.409
$ws.Cells.Item(113, 9).Value = 9713.333000000001  # I113: 509134.34 -> 9713.333000000001
$ws.Cells.Item(113, 11).Value = 9713.333000000001  # K113: 509134.34 -> 9713.333000000001
$ws.Cells.Item(113, 13).Value = -7543.333000000001  # M113: -506964.34 -> -7543.333000000001
$ws.Cells.Item(132, 8).Value = 6124.2085  # H132: 10025.4 -> 6124.2085
$ws.Cells.Item(132, 9).Value = 3298.8  # I132: 3583.3333 -> 3298.8
$ws.Cells.Item(132, 10).Value = 10833.223  # J132: 12786.286 -> 10833.223
$ws.Cells.Item(132, 11).Value = 9896.400000000001  # K132: 10749.9999 -> 9896.400000000001
$ws.Cells.Item(132, 12).Value = 32499.669  # L132: 38358.858 -> 32499.669
$ws.Cells.Item(132, 13).Value = -7366.400000000001  # M132: -8219.999899999999 -> -7366.400000000001
$ws.Cells.Item(132, 14).Value = -37559.669  # N132: -43418.858 -> -37559.669

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 0  # H5: 10000 -> 0
$ws.Cells.Item(5, 10).Value = 0  # J5: 10000 -> 0
$ws.Cells.Item(5, 12).Value = 0  # L5: 10000 -> 0
$ws.Cells.Item(5, 14).ClearContents()  # N5: -10224 -> (removed)
$ws.Cells.Item(6, 8).Value = 3750  # H6: 0 -> 3750
$ws.Cells.Item(6, 10).Value = 3750  # J6: 0 -> 3750
$ws.Cells.Item(6, 12).Value = 3750  # L6: 0 -> 3750
$ws.Cells.Item(6, 14).Value = -3980  # N6: None -> -3980
$ws.Cells.Item(132, 8).Value = 172325.52  # H132: 437565.4 -> 172325.52
$ws.Cells.Item(132, 9).Value = 191112.48  # I132: 490481.94 -> 191112.48
$ws.Cells.Item(132, 10).Value = 6373.8335  # J132: 3649.6 -> 6373.8335
$ws.Cells.Item(132, 11).Value = 573337.4400000001  # K132: 1471445.82 -> 573337.4400000001
$ws.Cells.Item(132, 12).Value = 19121.5005  # L132: 10948.8 -> 19121.5005
$ws.Cells.Item(132, 13).Value = -570807.4400000001  # M132: -1468915.82 -> -570807.4400000001
$ws.Cells.Item(132, 14).Value = -24181.5005  # N132: -16008.8 -> -24181.5005
$ws.Cells.Item(136, 8).Value = 2645.3396  # H136: 2994.7083 -> 2645.3396
$ws.Cells.Item(136, 9).Value = 1927.3721  # I136: 2350.139 -> 1927.3721
$ws.Cells.Item(136, 10).Value = 5732.6  # J136: 4928.4165 -> 5732.6
$ws.Cells.Item(136, 11).Value = 5782.1163  # K136: 7050.139 -> 5782.1163
$ws.Cells.Item(136, 12).Value = 17197.8  # L136: 14785.2495 -> 17197.8
